{"js": "// Update the stack-trace line numbers embedded in the big error-message\n// paragraph so they match the shifted line numbers in the regenerated\n// M2Doc sources (M2DocEvaluator / M2DocUtils / AbstractTemplatesTestSuite).\n//\n// Each pair below is unique within the document body, so a plain\n// text search-and-replace is safe (no risk of touching an unrelated\n// occurrence, even though \"M2DocEvaluator.java:1038)\" itself occurs\n// three times \u2014 every one of those three must become \":1084)\" anyway).\n\nconst replacements = [\n  [\"M2DocEvaluator.java:1489)\", \"M2DocEvaluator.java:1535)\"],\n  [\"M2DocEvaluator.java:1038)\", \"M2DocEvaluator.java:1084)\"],\n  [\"M2DocEvaluator.java:1254)\", \"M2DocEvaluator.java:1300)\"],\n  [\"M2DocEvaluator.java:275)\", \"M2DocEvaluator.java:278)\"],\n  [\"M2DocEvaluator.java:264)\", \"M2DocEvaluator.java:267)\"],\n  [\"M2DocUtils.java:712)\", \"M2DocUtils.java:694)\"],\n  [\"AbstractTemplatesTestSuite.java:459)\", \"AbstractTemplatesTestSuite.java:475)\"],\n  [\"AbstractTemplatesTestSuite.java:369)\", \"AbstractTemplatesTestSuite.java:384)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  // eslint-disable-next-line no-await-in-loop\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n}\n", "ps1": "# Update the stack-trace line numbers embedded in the big error-message\n# paragraph so they match the shifted line numbers in the regenerated\n# M2Doc sources (M2DocEvaluator / M2DocUtils / AbstractTemplatesTestSuite).\n#\n# Each pair below is unique within the document body, so a plain\n# text find-and-replace is safe (no risk of touching an unrelated\n# occurrence, even though \"M2DocEvaluator.java:1038)\" itself occurs\n# three times - every one of those three must become \":1084)\" anyway).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"M2DocEvaluator.java:1489)\", \"M2DocEvaluator.java:1535)\"),\n    @(\"M2DocEvaluator.java:1038)\", \"M2DocEvaluator.java:1084)\"),\n    @(\"M2DocEvaluator.java:1254)\", \"M2DocEvaluator.java:1300)\"),\n    @(\"M2DocEvaluator.java:275)\", \"M2DocEvaluator.java:278)\"),\n    @(\"M2DocEvaluator.java:264)\", \"M2DocEvaluator.java:267)\"),\n    @(\"M2DocUtils.java:712)\", \"M2DocUtils.java:694)\"),\n    @(\"AbstractTemplatesTestSuite.java:459)\", \"AbstractTemplatesTestSuite.java:475)\"),\n    @(\"AbstractTemplatesTestSuite.java:369)\", \"AbstractTemplatesTestSuite.java:384)\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
